# Scheduled-runner refresh of cached FFXIV market-board figures in the
# Zeromus_Profits leve-profit tables (currentAveragePrice* / LevePrice* /
# LeveProfit* columns, H:N) across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 28000
$ws.Range("J120").Value = 28000
$ws.Range("L120").Value = 28000
$ws.Range("N120").Value = -37676

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 857.1579
$ws.Range("I2").Value = 571.5333000000001
$ws.Range("J2").Value = 1928.25
$ws.Range("K2").Value = 571.5333000000001
$ws.Range("L2").Value = 1928.25
$ws.Range("M2").Value = -458.5333000000001
$ws.Range("N2").Value = -2154.25

$ws.Range("H107").Value = 37166.668
$ws.Range("J107").Value = 37166.668
$ws.Range("L107").Value = 37166.668
$ws.Range("N107").Value = -44846.668

$ws.Range("H108").Value = 55000
$ws.Range("J108").Value = 55000
$ws.Range("L108").Value = 55000
$ws.Range("N108").Value = -62680

$ws.Range("H109").Value = 34062.5
$ws.Range("J109").Value = 34062.5
$ws.Range("L109").Value = 34062.5
$ws.Range("N109").Value = -36836.5

$ws.Range("H110").Value = 2182.3
$ws.Range("I110").Value = 1741.1428
$ws.Range("J110").Value = 3211.6667
$ws.Range("K110").Value = 1741.1428
$ws.Range("L110").Value = 3211.6667
$ws.Range("M110").Value = 303.8571999999999
$ws.Range("N110").Value = -7301.6667

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H116").Value = 857.1579
$ws.Range("I116").Value = 571.5333000000001
$ws.Range("J116").Value = 1928.25
$ws.Range("K116").Value = 571.5333000000001
$ws.Range("L116").Value = 1928.25
$ws.Range("M116").Value = 1722.4667
$ws.Range("N116").Value = -6516.25

$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

$ws.Range("H119").Value = 48000
$ws.Range("J119").Value = 48000
$ws.Range("L119").Value = 48000
$ws.Range("N119").Value = -57676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 857.1579
$ws.Range("I3").Value = 571.5333000000001
$ws.Range("J3").Value = 1928.25
$ws.Range("K3").Value = 571.5333000000001
$ws.Range("L3").Value = 1928.25
$ws.Range("M3").Value = -457.5333000000001
$ws.Range("N3").Value = -2156.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 564.19446
$ws.Range("I122").Value = 407.8095
$ws.Range("J122").Value = 783.13336
$ws.Range("K122").Value = 3670.2855
$ws.Range("L122").Value = 7048.20024
$ws.Range("M122").Value = -1220.2855
$ws.Range("N122").Value = -11948.20024

$ws.Range("H131").Value = 1516186.6
$ws.Range("I131").Value = 7407837
$ws.Range("J131").Value = 1190.8
$ws.Range("K131").Value = 22223511
$ws.Range("L131").Value = 3572.4
$ws.Range("M131").Value = -22218471
$ws.Range("N131").Value = -13652.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 916.2105
$ws.Range("I107").Value = 783.4545000000001
$ws.Range("J107").Value = 1098.75
$ws.Range("K107").Value = 783.4545000000001
$ws.Range("L107").Value = 1098.75
$ws.Range("M107").Value = 1136.5455
$ws.Range("N107").Value = -4938.75

$ws.Range("H109").Value = 12690
$ws.Range("J109").Value = 12690
$ws.Range("L109").Value = 12690
$ws.Range("N109").Value = -14770

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H111").Value = 24293
$ws.Range("J111").Value = 24293
$ws.Range("L111").Value = 24293
$ws.Range("N111").Value = -30427

$ws.Range("H113").Value = 1623.6875
$ws.Range("I113").Value = 758.9231
$ws.Range("J113").Value = 5371
$ws.Range("K113").Value = 758.9231
$ws.Range("L113").Value = 5371
$ws.Range("M113").Value = 1411.0769
$ws.Range("N113").Value = -9711

$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678

$ws.Range("H116").Value = 50742
$ws.Range("J116").Value = 50742
$ws.Range("L116").Value = 50742
$ws.Range("N116").Value = -59920

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1540.826
$ws.Range("I61").Value = 1605.9048
$ws.Range("J61").Value = 857.5
$ws.Range("K61").Value = 1605.9048
$ws.Range("L61").Value = 857.5
$ws.Range("M61").Value = -1403.9048
$ws.Range("N61").Value = -1261.5

$ws.Range("H108").Value = 52500.5
$ws.Range("I108").Value = 20001
$ws.Range("J108").Value = 85000
$ws.Range("K108").Value = 20001
$ws.Range("L108").Value = 85000
$ws.Range("M108").Value = -16161
$ws.Range("N108").Value = -92680

$ws.Range("H110").Value = 26982.75
$ws.Range("J110").Value = 26982.75
$ws.Range("L110").Value = 26982.75
$ws.Range("N110").Value = -35162.75

$ws.Range("H111").Value = 30465.25
$ws.Range("J111").Value = 30465.25
$ws.Range("L111").Value = 30465.25
$ws.Range("N111").Value = -38645.25

$ws.Range("H112").Value = 34700
$ws.Range("J112").Value = 34700
$ws.Range("L112").Value = 34700
$ws.Range("N112").Value = -37654

$ws.Range("H113").Value = 1540.826
$ws.Range("I113").Value = 1605.9048
$ws.Range("J113").Value = 857.5
$ws.Range("K113").Value = 1605.9048
$ws.Range("L113").Value = 857.5
$ws.Range("M113").Value = 564.0952
$ws.Range("N113").Value = -5197.5

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H115").Value = 1000000
$ws.Range("J115").Value = 1000000
$ws.Range("L115").Value = 1000000
$ws.Range("N115").Value = -1002350

$ws.Range("H136").Value = 4451.4
$ws.Range("I136").Value = 8288.294
$ws.Range("J136").Value = 1615.4348
$ws.Range("K136").Value = 24864.882
$ws.Range("L136").Value = 4846.3044
$ws.Range("M136").Value = -22314.882
$ws.Range("N136").Value = -9946.304400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
